$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the multiplication table body (B2:I9).
# Row header (A2:A9) already contains 1..8, column header (B1:I1) already contains 1..8.
for ($row = 2; $row -le 9; $row++) {
    $rowMultiplier = $row - 1
    for ($col = 2; $col -le 9; $col++) {
        $colMultiplier = $col - 1
        $ws.Cells.Item($row, $col).Value = $rowMultiplier * $colMultiplier
    }
}
